# atualizações de funções de tratamento das colunas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Insert a new row at position 8 (pushing "pais", "email", "dob.date" rows down by one)
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the "rua" (street) column metadata
$ws.Cells.Item(8, 1).Value = "cadastro"
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 4).Value = "rua"
$ws.Cells.Item(8, 5).Value = "rua do endereço do indivíduo"
$ws.Cells.Item(8, 6).Value = "string"
$ws.Cells.Item(8, 7).Value = 0

# Re-number the "id" column for the rows that were pushed down
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(11, 2).Value = 10

# Add the new "limpa_texto" column (H) with its header
$ws.Cells.Item(1, 8).Value = "limpa_texto"
$ws.Cells.Item(1, 8).Font.Bold = $true
$ws.Cells.Item(1, 8).Font.Italic = $true
$ws.Cells.Item(1, 8).HorizontalAlignment = -4108

# "nome_original" for the new row is set last (affects shared-string ordering)
$ws.Cells.Item(8, 3).Value = "location.street.name"

# Fill the "limpa_texto" flag values for all data rows
for ($r = 2; $r -le 11; $r++) {
    if ($r -eq 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}

# Update the selected cell/view
$ws.Range("H8").Select()
